# AugerQuant_function_descriptions.xlsx
# Commit: "changes to handling of kwargs for scatter plots"
#
# A new row documenting the `get_plotkwargs` helper function is inserted
# into the "Functions" sheet at row 103 (pushing everything from the old
# row 104 onward down by one row). Row-anchored metadata (comments, the
# hidden _FilterDatabase defined name, the worksheet dimension, and the
# current selection) are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert the new row. This shifts all cell contents/styles below it
#    down by one (rows 104-163 -> 105-164) and updates <dimension> for us.
# ---------------------------------------------------------------------
$ws.Rows.Item(103).Insert()

# ---------------------------------------------------------------------
# 2) Populate the newly inserted row 103 with the new function entry.
#    Columns: A=Module, B=Functions, C=Calling funct, D=Args,
#             E=subfunctions, F=Returns, G=Status, H=Description, I=TODO
# ---------------------------------------------------------------------
$ws.Range("A103").Value = "AESplot"
$ws.Range("B103").Value = "get_plotkwargs"
$ws.Range("C103").Value = "scattercompplot, others?"
$ws.Range("H103").Value = "find and set x and y error columns (if requested by passed kwargs)"

# ---------------------------------------------------------------------
# 3) Move the comments that live on/after the old row 104 down by one
#    row so they stay attached to the same logical content (the row
#    insert above does not relocate cell comments automatically).
# ---------------------------------------------------------------------
$commentRefs = @("H106","G109","H113","H114","H116","H117","H139","H141","H146")

$movedComments = @()
foreach ($ref in $commentRefs) {
    if ($ref -match '^([A-Z]+)([0-9]+)$') {
        $col = $Matches[1]
        $rowNum = [int]$Matches[2]
        $newRef = "$col$($rowNum + 1)"
        $cmt = $ws.Range($ref).Comment
        $text = $cmt.Text()
        $movedComments += , @($newRef, $text)
        $cmt.Delete()
    }
}
foreach ($pair in $movedComments) {
    $ws.Range($pair[0]).AddComment($pair[1]) | Out-Null
}

# ---------------------------------------------------------------------
# 4) Update the hidden _xlnm._FilterDatabase defined name so its range
#    covers the new last row (163 -> 164).
# ---------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -like "*_FilterDatabase*") {
        $nm.RefersTo = "=Functions!`$A`$1:`$I`$164"
    }
}

# ---------------------------------------------------------------------
# 5) Restore the view's active selection (bottom-right pane) onto the
#    cell that now holds the shifted H104 -> H105 entry.
# ---------------------------------------------------------------------
$ws.Range("H105").Select()

Write-Host "Inserted get_plotkwargs row and shifted trailing rows/comments."
